$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename headers L1:N1 from stat_u/stat_plus/stat_minus to tot_u/tot_plus/tot_minus
$ws.Range("L1").Value = "tot_u"
$ws.Range("M1").Value = "tot_plus"
$ws.Range("N1").Value = "tot_minus"

# Delete column O (the %syst_c column) - shifts everything left
$ws.Range("O1").EntireColumn.Delete()

# Reset selection to N2
$ws.Range("N2").Select() | Out-Null
